$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[60.48182471966794, 66.42747218271614]"
$ws.Range("T2").Value = "[47.83340699538164, 51.94758269242783]"
$ws.Range("L3").Value = "[59.61900046690208, 68.28833901308474]"
$ws.Range("T3").Value = "[47.988484819676934, 52.7108142540737]"
